$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update Version value 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Update Date value
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicate "Contact" / "No display for ContactDetail" row (row 11)
$ws1.Rows.Item(11).Delete()

# Row 9 (Publisher) now gets a value; row 10 becomes Jurisdiction / United States of America
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"
